$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 359, shifting existing rows 359..469 down to 360..470.
$ws.Rows.Item(359).Insert()

# Seed the new row with the same "constant" columns as the surrounding data
# (every row in this sheet shares the same market/category/unit metadata).
for ($col = 1; $col -le 18; $col++) {
    $ws.Cells.Item(359, $col).Value = $ws.Cells.Item(358, $col).Value2
}

# Now overwrite the columns that actually carry new data for this record.
$ws.Cells.Item(359, 4).Value = 45093    # D359 Fecha
$ws.Cells.Item(359, 10).Value = 165     # J359 Volumen
$ws.Cells.Item(359, 11).Value = 7500    # K359 Precio minimo
$ws.Cells.Item(359, 12).Value = 8000    # L359 Precio maximo
$ws.Cells.Item(359, 13).Value = 7758    # M359 Precio promedio ponderado
$ws.Cells.Item(359, 16).Value = 129     # P359 Precio $/Kg
